$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row: D (Fecha serial), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), S (Precio $/Kg)
$data = @{
    2 = @(44232, 60, 11000, 12000, 11583, 827)
    3 = @(44216, 55, 11000, 12000, 11545, 825)
    4 = @(44229, 55, 11000, 12000, 11364, 812)
    5 = @(44253, 90, 12000, 13000, 12667, 905)
    6 = @(44172, 90, 8500, 9000, 8806, 629)
    7 = @(44181, 65, 9000, 10000, 9462, 676)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals[0]
    $ws.Range("M$row").Value = $vals[1]
    $ws.Range("N$row").Value = $vals[2]
    $ws.Range("O$row").Value = $vals[3]
    $ws.Range("P$row").Value = $vals[4]
    $ws.Range("S$row").Value = $vals[5]
}
